# "spring 2017 made current"
# Shift the status/lastmod/priority columns down one row: the previously
# "future" Spring_2017 row becomes "current", the previously "current"
# Fall_2016 row becomes "past" (and gets a lastmod date + updated priority),
# and the remaining "past" rows' priority values shift down too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Spring_2017 -> current, priority 0.6
$ws.Range("B2").Value = "current"
$ws.Range("D2").Value = 0.6

# Row 3: Fall_2016 -> past, lastmod 12/10/2016 (serial 42714), priority 0.3
# Copy the number format already used by the lastmod column (row 4) so the
# new date cell keeps the workbook's existing date style instead of Excel
# inventing a new one from a parsed date string.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B3").Value = "past"
$ws.Range("C3").Value = 42714
$ws.Range("D3").Value = 0.3

# Row 4: Spring_2016 priority -> 0.1
$ws.Range("D4").Value = 0.1

# Row 5: Fall_2015 priority -> 0.05
$ws.Range("D5").Value = 0.05

# Update the active selection left behind when the file was saved
$ws.Range("E7").Select()
